# Update the "想去人数" (want-to-go count) figures in both the "展览"
# sheet (1st worksheet) and the "全部类型" sheet (4th worksheet), which
# duplicates the same rows. Values below reflect a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 1154
$ws1.Range("F4").Value  = 263
$ws1.Range("F5").Value  = 145
$ws1.Range("F6").Value  = 6
$ws1.Range("F7").Value  = 12202
$ws1.Range("F10").Value = 123
$ws1.Range("F11").Value = 11994
$ws1.Range("F12").Value = 4795
$ws1.Range("F13").Value = 2601
$ws1.Range("F14").Value = 112
$ws1.Range("F15").Value = 43
$ws1.Range("F18").Value = 942
$ws1.Range("F20").Value = 160
$ws1.Range("F21").Value = 68

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 1154
$ws4.Range("F4").Value  = 263
$ws4.Range("F5").Value  = 145
$ws4.Range("F8").Value  = 6
$ws4.Range("F9").Value  = 12202
$ws4.Range("F12").Value = 123
$ws4.Range("F13").Value = 11994
$ws4.Range("F14").Value = 4795
$ws4.Range("F15").Value = 2603
$ws4.Range("F16").Value = 112
$ws4.Range("F17").Value = 43
$ws4.Range("F20").Value = 942
$ws4.Range("F22").Value = 160
$ws4.Range("F23").Value = 68
